$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks; they will be re-added after the table is
# rebuilt below, since row insertion does not shift hyperlink ranges in
# this engine.
$ws.Hyperlinks.Delete()

# Insert two new rows for the new accommodation entries (Visby Boers,
# Donners Hotell) above the existing "BW Strand Hotel" row, pushing the
# remaining rows (including Ferry terminal / Airport / Kapitelhusgaarden)
# down by two rows. Excel duplicates the formatting of the row above,
# which already matches column G's Hyperlink style used throughout the
# table.
$ws.Rows("3:4").Insert()

# Re-write every data row (2-15) with the final values so cell contents
# and ordering match the updated table exactly.
$ws.Range("A2").Value = "briefcase"
$ws.Range("B2").Value = "red"
$ws.Range("C2").Value = 57.6394006989332
$ws.Range("D2").Value = 18.2886166637887
$ws.Range("E2").Value = "Cramérgatan 3, 621 57 Visby"
$ws.Range("F2").Value = "Workshop Venue"
$ws.Range("G2").Value = "http://www.campusgotland.uu.se/"

$ws.Range("A3").Value = "bed"
$ws.Range("B3").Value = "blue"
$ws.Range("C3").Value = 57.6395457505725
$ws.Range("D3").Value = 18.2921123963061
$ws.Range("E3").Value = "Strandgatan 10, 621 56 Visby, 290m (4 min walk)"
$ws.Range("F3").Value = "Visby Börs"
$ws.Range("G3").Value = "https://visbybors.se/"

$ws.Range("A4").Value = "bed"
$ws.Range("B4").Value = "blue"
$ws.Range("C4").Value = 57.6388154392005
$ws.Range("D4").Value = 18.2915317571919
$ws.Range("E4").Value = "Donnersgatan 6, 621 57 Visby, 290m (4 min walk)"
$ws.Range("F4").Value = "Donners Hotell"
$ws.Range("G4").Value = "https://donnershotell.se/"

$ws.Range("A5").Value = "bed"
$ws.Range("B5").Value = "blue"
$ws.Range("C5").Value = 57.6418222910998
$ws.Range("D5").Value = 18.2924924826879
$ws.Range("E5").Value = "Strandgatan 34, 621 56 Visby (8 min walk, 2 min bike)"
$ws.Range("F5").Value = "BW Strand Hotel"
$ws.Range("G5").Value = "http://www.strandhotel.se/"

$ws.Range("A6").Value = "bed"
$ws.Range("B6").Value = "blue"
$ws.Range("C6").Value = 57.6386489348324
$ws.Range("D6").Value = 18.2909058341587
$ws.Range("E6").Value = "Strandgatan 6, 621 57 Visby, 260m (3 min walk)"
$ws.Range("F6").Value = "Clarion Hotel Visby"
$ws.Range("G6").Value = "https://www.strawberryhotels.com/hotels/sweden/visby/clarion-hotel-wisby/"

$ws.Range("A7").Value = "bed"
$ws.Range("B7").Value = "blue"
$ws.Range("C7").Value = 57.6431496715119
$ws.Range("D7").Value = 18.2959276710474
$ws.Range("E7").Value = "Smedjegatan 3, 621 55 Visby (11 min walk, 4 min bike)"
$ws.Range("F7").Value = "Hotel St. Clemens"
$ws.Range("G7").Value = "http://www.clemenshotell.se/"

$ws.Range("A8").Value = "bed"
$ws.Range("B8").Value = "blue"
$ws.Range("C8").Value = 57.6320705322625
$ws.Range("D8").Value = 18.28036866919
$ws.Range("E8").Value = "Färjeleden 3, 621 58 Visby, 950m (13 min walk, 3 min bike)"
$ws.Range("F8").Value = "Scandic Visby"
$ws.Range("G8").Value = "https://www.scandichotels.com/en/hotels/scandic-visby"

$ws.Range("A9").Value = "bed"
$ws.Range("B9").Value = "blue"
$ws.Range("C9").Value = 57.6364820625953
$ws.Range("D9").Value = 18.2894577665854
$ws.Range("E9").Value = "Visborgsgatan 1, 621 57 Visby, 450m (6 min walk)"
$ws.Range("F9").Value = "Hotell Slottsbacken"
$ws.Range("G9").Value = "http://hotellslottsbacken.se/"

$ws.Range("A10").Value = "bed"
$ws.Range("B10").Value = "blue"
$ws.Range("C10").Value = 57.6386907756162
$ws.Range("D10").Value = 18.2951541327037
$ws.Range("E10").Value = "Hästgatan 14, 621 56 Visby, 450m (6 min walk)"
$ws.Range("F10").Value = "Visby Logi & Vandrarhem"
$ws.Range("G10").Value = "http://www.visbylogi.se/"

$ws.Range("A11").Value = "bed"
$ws.Range("B11").Value = "blue"
$ws.Range("C11").Value = 57.6258703549948
$ws.Range("D11").Value = 18.2820583121987
$ws.Range("E11").Value = "Söderväg 56 A, 621 58 Visby, 1.7km (25 min walk, 8 min bike)"
$ws.Range("F11").Value = "First Hotel Kokoloko"
$ws.Range("G11").Value = "https://www.firsthotels.se/hotell/sverige/gotland/first-hotel-kokoloko/"

$ws.Range("A12").Value = "bed"
$ws.Range("B12").Value = "blue"
$ws.Range("C12").Value = 57.6554010828628
$ws.Range("D12").Value = 18.3079058678692
$ws.Range("E12").Value = "Snäckgärdsvägen 32, 621 55 Visby, 2.3km (31 min walk, 7 min bike)"
$ws.Range("F12").Value = "Visby Strandby"
$ws.Range("G12").Value = "http://www.visbystrandby.se/"

$ws.Range("A13").Value = "ship"
$ws.Range("B13").Value = "purple"
$ws.Range("C13").Value = 57.634235891617
$ws.Range("D13").Value = 18.2812459933587
$ws.Range("E13").Value = "Färjeleden 14, 621 58 Visby"
$ws.Range("F13").Value = "Ferry terminal"
$ws.Range("G13").Value = "http://www.destinationgotland.se/"

$ws.Range("A14").Value = "plane"
$ws.Range("B14").Value = "purple"
$ws.Range("C14").Value = 57.6607926854904
$ws.Range("D14").Value = 18.33859496282
$ws.Range("E14").Value = "621 41 Visby"
$ws.Range("F14").Value = "Airport"
$ws.Range("G14").Value = "https://www.swedavia.com/visby/"

$ws.Range("A15").Value = "cutlery"
$ws.Range("B15").Value = "green"
$ws.Range("C15").Value = 57.641785617709
$ws.Range("D15").Value = 18.2960642268324
$ws.Range("E15").Value = "S:t Drottensgatan 8, 621 56 Visby"
$ws.Range("F15").Value = "Kapitelhusgården"
$ws.Range("G15").Value = "http://www.kapitelhusgarden.se/"


# Re-create the hyperlinks on column G for every row that has a URL.
$ws.Hyperlinks.Add($ws.Range("G2"), "http://www.campusgotland.uu.se/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), "https://visbybors.se/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G4"), "https://donnershotell.se/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G5"), "http://www.strandhotel.se/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G6"), "https://www.strawberryhotels.com/hotels/sweden/visby/clarion-hotel-wisby/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G7"), "http://www.clemenshotell.se/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G8"), "https://www.scandichotels.com/en/hotels/scandic-visby") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G9"), "http://hotellslottsbacken.se/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G10"), "http://www.visbylogi.se/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G11"), "https://www.firsthotels.se/hotell/sverige/gotland/first-hotel-kokoloko/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G12"), "http://www.visbystrandby.se/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G13"), "http://www.destinationgotland.se/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G14"), "https://www.swedavia.com/visby/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G15"), "http://www.kapitelhusgarden.se/") | Out-Null

# In the original workbook, the Kapitelhusgaarden row (G15, previously
# G13) keeps the plain "Normal" look instead of the blue/underlined
# Hyperlink style used by every other row, even though it is a working
# hyperlink. Match that quirk.
$ws.Range("G15").Style = "Normal"
$ws.Range("G15").VerticalAlignment = -4108    # xlCenter
$ws.Range("G15").HorizontalAlignment = -4131  # xlLeft

# Update the "info" defined name so it covers the two new rows.
$wb.Names.Item("Sheet1!info").RefersTo = "=Sheet1!`$A`$1:`$F`$13"

# Update the active selection to match the authored state.
$ws.Range("F19:F20").Select() | Out-Null
